$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

# Uniform employee naming convention: "Employee 1".."Employee 6" -> "Employee A".."Employee E", "Employee F", "Employee G"
# (the target sheets Red/Green already use the Employee A-E convention)
$ws.Range("A2").Value = "Employee A"
$ws.Range("A3").Value = "Employee B"
$ws.Range("A4").Value = "Employee C"
$ws.Range("A5").Value = "Employee D"
$ws.Range("A6").Value = "Employee F"
$ws.Range("A7").Value = "Employee G"

# Move the active selection on the Personnel sheet from D13 to A8
$ws.Activate()
$ws.Range("A8").Select()
